$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 200, pushing existing rows 200-219 down to 201-220.
$ws.Rows.Item(200).Insert()

# Populate the newly inserted row 200 with the new record.
$ws.Cells.Item(200, 1).Value = 1
$ws.Cells.Item(200, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(200, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(200, 4).Value = 45106
$ws.Cells.Item(200, 5).Value = 15
$ws.Cells.Item(200, 6).Value = "Fruta"
$ws.Cells.Item(200, 7).Value = 100106
$ws.Cells.Item(200, 8).Value = "Oleaginosos"
$ws.Cells.Item(200, 9).Value = 100106002
$ws.Cells.Item(200, 10).Value = "Palta"
$ws.Cells.Item(200, 11).Value = "Hass"
$ws.Cells.Item(200, 12).Value = "Primera"
$ws.Cells.Item(200, 13).Value = 208
$ws.Cells.Item(200, 14).Value = 20000
$ws.Cells.Item(200, 15).Value = 22000
$ws.Cells.Item(200, 16).Value = 21000
$ws.Cells.Item(200, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(200, 18).Value = "Perú"
$ws.Cells.Item(200, 19).Value = 2100
$ws.Cells.Item(200, 20).Value = 10
